$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.13508
$ws.Range("H2").Value = 3.40524
$ws.Range("I2").Value = 0.1224366388308639
$ws.Range("J2").Value = 0.1224366388308639
$ws.Range("M2").Value = 2.633709
$ws.Range("N2").Value = 7.901127
$ws.Range("O2").Value = 0.3696831990833914
$ws.Range("P2").Value = 0.3696831990833913
$ws.Range("Q2").Value = 2.98947041172
$ws.Range("R2").Value = 26.90523370548
$ws.Range("S2").Value = 0.04526276832801156
$ws.Range("T2").Value = 0.04526276832801155

$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.13508
$ws.Range("H3").Value = 3.40524
$ws.Range("I3").Value = 0.1224366388308639
$ws.Range("J3").Value = 0.1224366388308639
$ws.Range("M3").Value = 4.488244
$ws.Range("N3").Value = 13.464732
$ws.Range("O3").Value = 0.6299968600125665
$ws.Range("P3").Value = 0.6299968600125665
$ws.Range("Q3").Value = 5.09451599952
$ws.Range("R3").Value = 45.85064399568
$ws.Range("S3").Value = 0.07713469801393696
$ws.Range("T3").Value = 0.07713469801393695

$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.13508
$ws.Range("H4").Value = 3.40524
$ws.Range("I4").Value = 0.1224366388308639
$ws.Range("J4").Value = 0.1224366388308639
$ws.Range("N4").Value = 0.006838
$ws.Range("O4").Value = 0.0003199409040421993
$ws.Range("P4").Value = 0.0003199409040421993
$ws.Range("Q4").Value = 0.00258722568
$ws.Range("R4").Value = 0.02328503112
$ws.Range("S4").Value = 0.00003917248891543485
$ws.Range("T4").Value = 0.00003917248891543485

$ws.Range("G5").Value = 0.672624
$ws.Range("H5").Value = 2.017872
$ws.Range("I5").Value = 0.07255331937570129
$ws.Range("J5").Value = 0.07255331937570129
$ws.Range("M5").Value = 2.633709
$ws.Range("N5").Value = 7.901127
$ws.Range("O5").Value = 0.3696831990833914
$ws.Range("P5").Value = 0.3696831990833913
$ws.Range("Q5").Value = 1.771495882416
$ws.Range("R5").Value = 15.943462941744
$ws.Range("S5").Value = 0.02682174321092826
$ws.Range("T5").Value = 0.02682174321092825

$ws.Range("G6").Value = 0.672624
$ws.Range("H6").Value = 2.017872
$ws.Range("I6").Value = 0.07255331937570129
$ws.Range("J6").Value = 0.07255331937570129
$ws.Range("M6").Value = 4.488244
$ws.Range("N6").Value = 13.464732
$ws.Range("O6").Value = 0.6299968600125665
$ws.Range("P6").Value = 0.6299968600125665
$ws.Range("Q6").Value = 3.018900632256
$ws.Range("R6").Value = 27.170105690304
$ws.Range("S6").Value = 0.04570836339018072
$ws.Range("T6").Value = 0.04570836339018072

$ws.Range("G7").Value = 0.672624
$ws.Range("H7").Value = 2.017872
$ws.Range("I7").Value = 0.07255331937570129
$ws.Range("J7").Value = 0.07255331937570129
$ws.Range("N7").Value = 0.006838
$ws.Range("O7").Value = 0.0003199409040421993
$ws.Range("P7").Value = 0.0003199409040421993
$ws.Range("Q7").Value = 0.001533134304
$ws.Range("R7").Value = 0.013798208736
$ws.Range("S7").Value = 0.00002321277459232429
$ws.Range("T7").Value = 0.00002321277459232429

$ws.Range("G8").Value = 7.463050333333332
$ws.Range("H8").Value = 22.389151
$ws.Range("I8").Value = 0.8050100417934347
$ws.Range("J8").Value = 0.8050100417934348
$ws.Range("M8").Value = 2.633709
$ws.Range("N8").Value = 7.901127
$ws.Range("O8").Value = 0.3696831990833914
$ws.Range("P8").Value = 0.3696831990833913
$ws.Range("Q8").Value = 19.655502830353
$ws.Range("R8").Value = 176.899525473177
$ws.Range("S8").Value = 0.2975986875444515
$ws.Range("T8").Value = 0.2975986875444515

$ws.Range("G9").Value = 7.463050333333332
$ws.Range("H9").Value = 22.389151
$ws.Range("I9").Value = 0.8050100417934347
$ws.Range("J9").Value = 0.8050100417934348
$ws.Range("M9").Value = 4.488244
$ws.Range("N9").Value = 13.464732
$ws.Range("O9").Value = 0.6299968600125665
$ws.Range("P9").Value = 0.6299968600125665
$ws.Range("Q9").Value = 33.49599088028133
$ws.Range("R9").Value = 301.4639179225319
$ws.Range("S9").Value = 0.5071537986084488
$ws.Range("T9").Value = 0.5071537986084489

$ws.Range("G10").Value = 7.463050333333332
$ws.Range("H10").Value = 22.389151
$ws.Range("I10").Value = 0.8050100417934347
$ws.Range("J10").Value = 0.8050100417934348
$ws.Range("N10").Value = 0.006838
$ws.Range("O10").Value = 0.0003199409040421993
$ws.Range("P10").Value = 0.0003199409040421993
$ws.Range("Q10").Value = 0.01701077939311111
$ws.Range("R10").Value = 0.153097014538
$ws.Range("S10").Value = 0.0002575556405344401
$ws.Range("T10").Value = 0.0002575556405344402
